# Update countries & provincias Spain
# Applies the 26-Jul-2020 13:14 -> 14:31 COVID data refresh to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp ------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 26 de Julio de 2020 a las 14:31"

# --- Country name swaps (ranking changed order) -----------------------
# Estado de Palestina overtook Bulgaria in total cases.
$ws.Range("A80").Value = "Estado de Palestina"
$ws.Range("A81").Value = "Bulgaria"

# Groenlandia / Islas Malvinas swapped rank position (tied totals).
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Islas Malvinas"

# --- Numeric data refresh ----------------------------------------------
# Row 4: Estados Unidos
$ws.Range("B4").Value = 4315926
$ws.Range("C4").Value = 217
$ws.Range("E4").Value = 2104834
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 149400

# Row 6: India
$ws.Range("B6").Value = 1396304
$ws.Range("C6").Value = 10810
$ws.Range("D6").Value = 892727
$ws.Range("E6").Value = 471370
$ws.Range("G6").Value = 111
$ws.Range("H6").Value = 32207

# Row 14: Iran
$ws.Range("B14").Value = 291172
$ws.Range("C14").Value = 2333
$ws.Range("D14").Value = 253213
$ws.Range("E14").Value = 22259
$ws.Range("G14").Value = 216
$ws.Range("H14").Value = 15700

# Row 55: Suiza
$ws.Range("B55").Value = 34412
$ws.Range("C55").Value = 110
$ws.Range("D55").Value = 30700
$ws.Range("E55").Value = 1735

# Row 65: Uzbekistan
$ws.Range("E65").Value = 9281
$ws.Range("G65").Value = 3
$ws.Range("H65").Value = 114

# Row 67: Nepal
$ws.Range("B67").Value = 18613
$ws.Range("C67").Value = 130
$ws.Range("D67").Value = 13128
$ws.Range("E67").Value = 5440

# Row 80: Estado de Palestina (new data for this rank)
$ws.Range("B80").Value = 10469
$ws.Range("C80").Value = 163
$ws.Range("D80").Value = 3752
$ws.Range("E80").Value = 6642
$ws.Range("H80").Value = 75

# Row 81: Bulgaria (new data for this rank)
$ws.Range("B81").Value = 10312
$ws.Range("D81").Value = 5306
$ws.Range("E81").Value = 4668
$ws.Range("H81").Value = 338

# Row 84: Senegal
$ws.Range("B84").Value = 9681
$ws.Range("C84").Value = 129
$ws.Range("D84").Value = 6409
$ws.Range("E84").Value = 3081
$ws.Range("G84").Value = 4
$ws.Range("H84").Value = 191

# Row 88: Consejo Danes para los Refugiados
$ws.Range("B88").Value = 8831
$ws.Range("C88").Value = 30
$ws.Range("D88").Value = 5510
$ws.Range("E88").Value = 3117

# Row 98: Croacia
$ws.Range("B98").Value = 4857
$ws.Range("C98").Value = 65
$ws.Range("D98").Value = 3866
$ws.Range("E98").Value = 855
$ws.Range("G98").Value = 3
$ws.Range("H98").Value = 136

# Row 125: Eslovenia
$ws.Range("B125").Value = 2082
$ws.Range("C125").Value = 16
$ws.Range("E125").Value = 288
